$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.803.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '''1.635.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''215.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("D10").Value = '''19.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("D11").Value = '''0.0793'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").Value = '''4.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '''1.636.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").Value = '''1.860.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '''62.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '''25.819.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '''4.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("D21").Value = '''194.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '''9.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''6.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.46%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("E25").Value = '  +2.08%  '
$ws.Range("D26").Value = '''142.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.21%  '
$ws.Range("D27").Value = '''0.125'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").Value = '''6.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").Value = '''15.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("D32").Value = '''3.37'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.94%  '
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("D37").Value = '''1.137.93'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("E40").Value = '  -0.42%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  +1.90%  '
$ws.Range("D43").Value = '''100.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").Value = '''0.808'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").Value = '''1.770.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").Value = '''0.0₆0112'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.99%  '
$ws.Range("D47").Value = '''55.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("E49").Value = '  -2.01%  '
$ws.Range("D50").Value = '''7.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.80%  '
$ws.Range("E51").Value = '  +0.03%  '
